$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
# D-column price values are plain text in the source data (inlineStr);
# force text format so numeric-looking strings (e.g. "1.003", "1.000")
# are not silently coerced into numbers, which would lose formatting.
$ws.Range("D2").Value = "28.152.21"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "1.797.78"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.10"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4292"
$ws.Range("E7").Value = "  -3.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3635"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.68"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07554"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.120"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.74"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.176"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.368"
$ws.Range("D16").Value = "1.813.56"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.69"
$ws.Range("E17").Value = "  +3.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001072"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06350"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.25"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.996"
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("D23").Value = "28.183.92"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.41"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("E25").Value = "  -5.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.40"
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.39"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").Value = "2.017.54"
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.237"
$ws.Range("E29").Value = "  -5.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.83"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.178"
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.868"
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09030"
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.550"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.79"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02361"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.125"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6509"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2123"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06119"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.193"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.430"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.949"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.66"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6031"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.711"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.57"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.001"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.155"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06972"
$ws.Range("E51").Value = "  +0.98%  "
